# Allegato D — aggiornamento calcolo dei garanti
# Applies the corrections described in the commit diff:
#  - fixes a handful of mistyped CLASSE codes
#  - normalises "scientifico tecnologica" -> "scientifico tecnologico"
#    for the LM rows that were missing it
#  - fills in the previously-empty "N. di riferimento" / "N. max"
#    columns (F/G) for those same rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple CLASSE (column C) text corrections -----------------------
$ws.Range("C5").Value  = "L/SNT4"
$ws.Range("C17").Value = "L-Sc.Mat."
$ws.Range("C26").Value = "L-4"
$ws.Range("C52").Value = "L-P01"
$ws.Range("C53").Value = "L-P02"
$ws.Range("C54").Value = "L-P03"
$ws.Range("C59").Value = "LM-17"
$ws.Range("C157").Value = "LM-42"

# --- Column A: "scientifico tecnologica" -> "scientifico tecnologico" -
# and fill in columns F (N. di riferimento) / G (N. max) as text, since
# every other value in those columns is stored as text in this sheet.

for ($r = 59; $r -le 87; $r++) {
    $ws.Cells.Item($r, 1).Value = "scientifico tecnologico"
    $ws.Cells.Item($r, 6).NumberFormat = "@"
    $ws.Cells.Item($r, 6).Value = "65"
    $ws.Cells.Item($r, 7).NumberFormat = "@"
    $ws.Cells.Item($r, 7).Value = "65"
}

for ($r = 88; $r -le 113; $r++) {
    $ws.Cells.Item($r, 1).Value = "scientifico tecnologico"
    $ws.Cells.Item($r, 6).NumberFormat = "@"
    $ws.Cells.Item($r, 6).Value = "65"
    $ws.Cells.Item($r, 7).NumberFormat = "@"
    $ws.Cells.Item($r, 7).Value = "80"
}

for ($r = 158; $r -le 160; $r++) {
    $ws.Cells.Item($r, 1).Value = "scientifico tecnologico"
    $ws.Cells.Item($r, 6).NumberFormat = "@"
    $ws.Cells.Item($r, 6).Value = "75"
    $ws.Cells.Item($r, 7).NumberFormat = "@"
    $ws.Cells.Item($r, 7).Value = "100"
}
